$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-10-28 Tuesday" "2025-10-29 Wednesday"

Replace-Text "543×2=1086" "104×7=728"
Replace-Text "102×2=204" "402×9=3618"
Replace-Text "628×5=3140" "455×6=2730"
Replace-Text "322×9=2898" "220×4=880"
Replace-Text "408×6=2448" "540×8=4320"

Replace-Text "692×5=3460" "589×9=5301"
Replace-Text "913×7=6391" "992×4=3968"
Replace-Text "788×7=5516" "893×7=6251"
Replace-Text "831×6=4986" "703×5=3515"
Replace-Text "342×6=2052" "713×8=5704"

Replace-Text "304×5=1520" "894×2=1788"
Replace-Text "963×7=6741" "472×5=2360"
Replace-Text "125×7=875" "909×2=1818"
Replace-Text "269×5=1345" "950×6=5700"
Replace-Text "947×4=3788" "445×8=3560"

Replace-Text "197×8=1576" "494×3=1482"
Replace-Text "189×8=1512" "519×4=2076"
Replace-Text "609×7=4263" "469×9=4221"
Replace-Text "774×9=6966" "141×7=987"
Replace-Text "942×5=4710" "254×9=2286"

Replace-Text "404×6=2424" "433×8=3464"
Replace-Text "398×9=3582" "162×5=810"
Replace-Text "697×6=4182" "700×3=2100"
Replace-Text "173×2=346" "535×4=2140"
Replace-Text "869×5=4345" "978×3=2934"

Write-Output "Replacements complete"
